$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'orch_session_id'
$ws.Range("D1").Value = 'orch_session_entry_id'
$ws.Range("H1").Value = 'orch_session_issue_id'
$ws.Range("H6").Value = 'db54c94e-e025-4555-92df-3b0472d624b8'
$ws.Range("H7").Value = '4c2d1bf8-0c52-4828-96ef-bf3f69d46d3d'
$ws.Range("H8").Value = 'bd7aab99-46b7-449b-9155-b153985b713a'
$ws.Range("H9").Value = '93fe09ef-5ae8-4848-91ce-9fb9ab3d850b'
$ws.Range("H10").Value = '671d7155-1e1d-4126-a5eb-874ebc5d57e4'
$ws.Range("H11").Value = 'baf75e04-d1ba-4794-bf82-467be42dcc75'
$ws.Range("H12").Value = '4c9c7ddf-1e6a-43a6-8c51-3fc78c145b00'
$ws.Range("H13").Value = '5a075599-6096-4961-8f04-96d8ef636f28'
$ws.Range("H14").Value = '0a763e51-1479-429d-816c-2e3217ca63f0'
$ws.Range("H15").Value = '5032c177-7c52-47de-8741-7d88464bc6fb'
$ws.Range("H16").Value = 'd1626583-bba0-4ec9-9fab-8d3d2c87c1d1'
$ws.Range("H17").Value = '778b95c3-3fc5-4c42-ac1f-65ac9ffb337f'
$ws.Range("H19").Value = '6116a37a-4e94-434c-be8d-fba574a57bce'
$ws.Range("H20").Value = 'bb80d8cd-9031-4498-a11b-38ae517e03f9'
$ws.Range("H21").Value = '57d5489c-47c5-4b5e-9d00-fc041385ae69'
$ws.Range("H22").Value = '62687b89-afce-4bd1-9e6c-62223f7cda02'
$ws.Range("H23").Value = '9af44f4a-6258-4de2-97fb-ea9d8177f53a'
$ws.Range("H24").Value = 'b3fc3589-7c0c-4771-a0d5-ad279731ed07'
$ws.Range("H25").Value = '88af3db6-67b0-4be4-b33f-994457cae455'
$ws.Range("H26").Value = '517e26e8-cc26-4f27-a93e-f56c69125462'
$ws.Range("H27").Value = '12bbdabc-2df6-4bb7-ae05-ce78637576d6'
$ws.Range("H28").Value = 'd6de092b-1701-4af9-87b4-4071f1b58824'
$ws.Range("H29").Value = '30dad59d-3b76-4375-9c0d-6497f09395ca'
$ws.Range("H30").Value = '3317574e-4eb9-41df-96d8-77f3f1d6c9c7'
$ws.Range("H31").Value = '8b2d7596-295c-408b-8271-26af2725bd3c'
$ws.Range("H32").Value = '032d1455-a432-4455-866d-636379cc0113'
$ws.Range("H33").Value = '2421138d-c86c-40f8-8e39-0554a6102f76'
$ws.Range("H34").Value = 'a58ce8f9-34a1-4af0-9a4a-d3a5529d0077'
$ws.Range("H35").Value = '057d87d0-a510-46e2-8f78-9e338600eee8'
$ws.Range("H36").Value = 'cafc25aa-fc22-46ac-891b-499239ff5185'
$ws.Range("H37").Value = '50c4fd24-d066-4db3-af9b-edbf4b3f0492'
$ws.Range("H38").Value = '7c58c117-1e20-4bfc-a928-5dfc96c7defe'
$ws.Range("H39").Value = '69912e9f-389d-4309-a46a-108e5cf22159'
$ws.Range("H40").Value = 'fd70db33-83de-4564-81a8-22016e1fb21b'
$ws.Range("H41").Value = 'ac4b863f-7c12-4bc0-8385-a953e1ca6c22'
$ws.Range("H42").Value = 'c3ef3362-8bc0-427f-a20c-e74abf3f398e'
$ws.Range("K42").Value = 'orch_session_entry_id'
$ws.Range("H43").Value = '165159a9-2592-401c-bd36-0621c4cc8b6d'
$ws.Range("H44").Value = 'fb89466b-6735-4a0d-a141-0e490d798c7b'
$ws.Range("H45").Value = 'cf44dd80-6825-4395-8237-be93d4e22374'
$ws.Range("H46").Value = '40cb2fbe-db4e-4e00-9846-b7b9ae3d3a49'
$ws.Range("K46").Value = 'orch_session_issue_id'
$ws.Range("H47").Value = '23108655-c46e-4507-aa83-96cd25985e8d'
$ws.Range("H48").Value = 'e89539f2-c751-4cf3-962c-c0cba00848b4'
$ws.Range("H49").Value = '06f6ffd3-ef5c-409f-b990-137697918f42'
$ws.Range("H50").Value = '70652dd0-f325-4b7a-981b-469e2453e8c1'
$ws.Range("H51").Value = '81650c85-fd50-471e-8fc8-e9a15156448a'
$ws.Range("H52").Value = '882d4a4b-6e7e-4822-878f-30a769069d3d'
$ws.Range("H53").Value = '94dd2556-9e09-4f07-96df-401430c7ecd9'
$ws.Range("H54").Value = 'e99b5a98-6e87-4d94-a7b7-0eedc5b71a77'
$ws.Range("H55").Value = '97c77468-b41c-4e4f-98bf-114f6a221b0b'
$ws.Range("H56").Value = '6a14248a-eefb-429d-959c-e0ccd57c6ae2'
$ws.Range("H57").Value = '50fda451-475b-448b-8b87-e70843e852f8'
$ws.Range("H58").Value = 'cf5a1b31-3567-4096-adc7-769e42e206e1'
$ws.Range("H59").Value = '87ec150d-1cd3-476e-8cb4-38ca49c550a9'
$ws.Range("H60").Value = '7772fe11-edb7-4004-8f15-7b6ada8d6282'
$ws.Range("H61").Value = '61beca90-96cb-431f-ab58-d0122791fd81'
$ws.Range("H62").Value = 'b4d48c0f-1002-4492-ab15-94459a117702'
$ws.Range("H63").Value = '601ad1d2-31e2-405d-8e4d-3c16008308a7'
$ws.Range("H64").Value = '39b53120-1f5a-426b-b2a9-8533089e09c7'
$ws.Range("H65").Value = 'aa6fce67-0d05-422a-af60-fce41e223634'
$ws.Range("H66").Value = 'f59499c5-ddf1-451e-919a-0deae4bb1665'
